# Disable feedbacks via control settings
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BAEPAbCiPC")

# electricity (row 2), heat (row 15), hydrogen (row 22) -> disable (0)
$ws.Range("B2").Value = 0
$ws.Range("B15").Value = 0
$ws.Range("B22").Value = 0
